$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column header "zero_age_depth" in F1
$ws.Range("F1").Value = "zero_age_depth"

# Set the width of column F to (best-fit) match the width from the diff (XML width 14.33203125).
# The host's ColumnWidth -> stored XML width mapping adds a fixed 5/6 character padding,
# so we back it out to land as close as possible on the target stored width.
$ws.Range("F1").EntireColumn.ColumnWidth = 13.498697916666666

# Move the active selection to F2, matching the post-edit selection in the diff
$ws.Range("F2").Select()
